# Auto-generated Excel COM-interop script applying numeric updates
# produced by a scheduled price-refresh run across the FFXIV crafting-leve sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 1638.7
$ws.Range("I58").Value = 1453.75
$ws.Range("J58").Value = 1762
$ws.Range("K58").Value = 4361.25
$ws.Range("L58").Value = 5286
$ws.Range("M58").Value = -4211.25
$ws.Range("N58").Value = -5586

$ws.Range("H97").Value = 1618.75
$ws.Range("J97").Value = 1618.75
$ws.Range("L97").Value = 4856.25
$ws.Range("N97").Value = -5848.25

$ws.Range("H103").Value = 250
$ws.Range("J103").Value = 250
$ws.Range("L103").Value = 750
$ws.Range("N103").Value = -1922

$ws.Range("H131").Value = 2197.25
$ws.Range("I131").Value = 2197.25
$ws.Range("K131").Value = 6591.75
$ws.Range("M131").Value = -1551.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4653.8335
$ws.Range("I32").Value = 4653.8335
$ws.Range("K32").Value = 4653.8335
$ws.Range("M32").Value = -4366.8335

$ws.Range("H45").Value = 2850.4285
$ws.Range("I45").Value = 3153.6667
$ws.Range("K45").Value = 3153.6667
$ws.Range("M45").Value = -2776.6667

$ws.Range("H97").Value = 999.8570999999999
$ws.Range("I97").Value = 874.8333
$ws.Range("K97").Value = 874.8333
$ws.Range("M97").Value = -378.8333

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3697
$ws.Range("I20").Value = 2796.8
$ws.Range("J20").Value = 5947.5
$ws.Range("K20").Value = 2796.8
$ws.Range("L20").Value = 5947.5
$ws.Range("M20").Value = -2549.8
$ws.Range("N20").Value = -6441.5

$ws.Range("H75").Value = 6056.857
$ws.Range("I75").Value = 6056.857
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 6056.857
$ws.Range("M75").Value = -5120.857
$ws.Range("N75").ClearContents()

$ws.Range("H78").Value = 6056.857
$ws.Range("I78").Value = 6056.857
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 18170.571
$ws.Range("L78").Value = 0
$ws.Range("M78").Value = -13490.571
$ws.Range("N78").ClearContents()

$ws.Range("H94").Value = 699.0909
$ws.Range("I94").Value = 745.5
$ws.Range("J94").Value = 235
$ws.Range("K94").Value = 745.5
$ws.Range("L94").Value = 235
$ws.Range("M94").Value = -294.5
$ws.Range("N94").Value = -1137

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 9020.75
$ws.Range("I25").Value = 10166.571
$ws.Range("J25").Value = 1000
$ws.Range("K25").Value = 10166.571
$ws.Range("L25").ClearContents()
$ws.Range("M25").Value = -9992.571
$ws.Range("N25").Value = -1348

$ws.Range("H31").Value = 2480.2307
$ws.Range("J31").Value = 4108.75
$ws.Range("L31").Value = 4108.75
$ws.Range("N31").Value = -4698.75

$ws.Range("H34").Value = 2480.2307
$ws.Range("J34").Value = 4108.75
$ws.Range("L34").Value = 4108.75
$ws.Range("N34").Value = -4512.75

$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()

$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()

$ws.Range("H99").Value = 5513.75
$ws.Range("I99").Value = 5057
$ws.Range("K99").Value = 5057
$ws.Range("M99").Value = -3559

$ws.Range("H124").Value = 107999.4
$ws.Range("J124").Value = 110499.5
$ws.Range("L124").Value = 110499.5
$ws.Range("N124").Value = -115409.5

$ws.Range("H126").Value = 5513.75
$ws.Range("I126").Value = 5057
$ws.Range("K126").Value = 15171
$ws.Range("M126").Value = -12701

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 413.9091
$ws.Range("I4").Value = 61.615383
$ws.Range("J4").Value = 922.7778
$ws.Range("K4").Value = 184.846149
$ws.Range("L4").Value = 2768.3334
$ws.Range("M4").Value = -72.846149
$ws.Range("N4").Value = -2992.3334

$ws.Range("H114").Value = 1000
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 937541.1
$ws.Range("I14").Value = 1208379.9
$ws.Range("J14").Value = 125025
$ws.Range("K14").Value = 1208379.9
$ws.Range("L14").Value = 125025
$ws.Range("M14").Value = -1208211.9
$ws.Range("N14").Value = -125361

$ws.Range("H70").Value = 8999
$ws.Range("I70").Value = 8999
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 8999
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -8729
$ws.Range("N70").ClearContents()

$ws.Range("H73").Value = 8999
$ws.Range("I73").Value = 8999
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 8999
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -8063
$ws.Range("N73").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1750
$ws.Range("I7").Value = 1750
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 1750
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -1638
$ws.Range("N7").ClearContents()

$ws.Range("H100").Value = 1247
$ws.Range("I100").Value = 1247
$ws.Range("K100").Value = 1247
$ws.Range("M100").Value = -706

$ws.Range("H126").Value = 1750
$ws.Range("I126").Value = 1750
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 5250
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -2780
$ws.Range("N126").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1606
$ws.Range("I96").Value = 960
$ws.Range("J96").Value = 2252
$ws.Range("K96").Value = 960
$ws.Range("L96").Value = 2252
$ws.Range("M96").Value = 413
$ws.Range("N96").Value = -4998

$ws.Range("H100").Value = 298
$ws.Range("I100").Value = 298
$ws.Range("K100").Value = 596
$ws.Range("M100").Value = -55

$ws.Range("H126").Value = 1595.3889
$ws.Range("I126").Value = 1416.1333
$ws.Range("K126").Value = 4248.3999
$ws.Range("M126").Value = -1778.3999
